$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format column E (rows 2-23) as Text so the numeric-looking strings
# are stored as shared-string text (t="s") rather than being auto-converted
# to numbers, matching the original built_in_total column semantics.
$eRange = $ws.Range("E2:E23")
$eRange.NumberFormat = "@"

$ws.Range("A2").Value = "87811004_0222_RO"
$ws.Range("B2").Value = 447
$ws.Range("C2").Value = "RON"
$ws.Range("D2").Value = 10192.91
$ws.Range("E2").Value = "10192.91"

$ws.Range("A3").Value = "87811004_0222_PE"
$ws.Range("B3").Value = 17
$ws.Range("C3").Value = "PEN"
$ws.Range("D3").Value = 169.05
$ws.Range("E3").Value = "169.05"

$ws.Range("A4").Value = "87811004_0222_HU"
$ws.Range("B4").Value = 1032
$ws.Range("C4").Value = "HUF"
$ws.Range("D4").Value = 2197785
$ws.Range("E4").Value = "2197785"

$ws.Range("A5").Value = "87811004_0222_EU"
$ws.Range("B5").Value = 790
$ws.Range("C5").Value = "EUR"
$ws.Range("D5").Value = 4280.03
$ws.Range("E5").Value = "4280.03"

$ws.Range("A6").Value = "87811004_0222_MX"
$ws.Range("B6").Value = 88
$ws.Range("C6").Value = "MXN"
$ws.Range("D6").Value = 5686.8
$ws.Range("E6").Value = "5686.8"

$ws.Range("A7").Value = "87811004_0222_LL"
$ws.Range("B7").Value = 38
$ws.Range("C7").Value = "USD"
$ws.Range("D7").Value = 100.1
$ws.Range("E7").Value = "100.1"

$ws.Range("A8").Value = "87811004_0222_BG"
$ws.Range("B8").Value = 7
$ws.Range("C8").Value = "BGN"
$ws.Range("D8").Value = 12.18
$ws.Range("E8").Value = "12.18"

$ws.Range("A9").Value = "87811004_0222_BR"
$ws.Range("B9").Value = 41
$ws.Range("C9").Value = "BRL"
$ws.Range("D9").Value = 678.02
$ws.Range("E9").Value = "678.02"

$ws.Range("A10").Value = "87811004_0222_CA"
$ws.Range("B10").Value = 333
$ws.Range("C10").Value = "CAD"
$ws.Range("D10").Value = 1720.6
$ws.Range("E10").Value = "1720.6"

$ws.Range("A11").Value = "87811004_0222_CZ"
$ws.Range("B11").Value = 17
$ws.Range("C11").Value = "CZK"
$ws.Range("D11").Value = 1162
$ws.Range("E11").Value = "1162"

$ws.Range("A12").Value = "87811004_0222_CL"
$ws.Range("B12").Value = 29
$ws.Range("C12").Value = "CLP"
$ws.Range("D12").Value = 43805
$ws.Range("E12").Value = "43805"

$ws.Range("A13").Value = "87811004_0222_CO"
$ws.Range("B13").Value = 24
$ws.Range("C13").Value = "COP"
$ws.Range("D13").Value = 183120
$ws.Range("E13").Value = "183120"

$ws.Range("A14").Value = "87811004_0222_NZ"
$ws.Range("B14").Value = 58
$ws.Range("C14").Value = "NZD"
$ws.Range("D14").Value = 217.42
$ws.Range("E14").Value = "217.42"

$ws.Range("A15").Value = "87811004_0222_AU"
$ws.Range("B15").Value = 357
$ws.Range("C15").Value = "AUD"
$ws.Range("D15").Value = 2286.56
$ws.Range("E15").Value = "2286.56"

$ws.Range("A16").Value = "87811004_0222_CH"
$ws.Range("B16").Value = 67
$ws.Range("C16").Value = "CHF"
$ws.Range("D16").Value = 281.01
$ws.Range("E16").Value = "281.01"

$ws.Range("A17").Value = "87811004_0222_NO"
$ws.Range("B17").Value = 30
$ws.Range("C17").Value = "NOK"
$ws.Range("D17").Value = 887.6
$ws.Range("E17").Value = "887.6"

$ws.Range("A18").Value = "87811004_0222_US"
$ws.Range("B18").Value = 1575
$ws.Range("C18").Value = "USD"
$ws.Range("D18").Value = 11197.9
$ws.Range("E18").Value = "11197.9"

$ws.Range("A19").Value = "87811004_0222_DK"
$ws.Range("B19").Value = 31
$ws.Range("C19").Value = "DKK"
$ws.Range("D19").Value = 478.8
$ws.Range("E19").Value = "478.8"

$ws.Range("A20").Value = "87811004_0222_PL"
$ws.Range("B20").Value = 44
$ws.Range("C20").Value = "PLN"
$ws.Range("D20").Value = 498.12
$ws.Range("E20").Value = "498.12"

$ws.Range("A21").Value = "87811004_0222_SE"
$ws.Range("B21").Value = 54
$ws.Range("C21").Value = "SEK"
$ws.Range("D21").Value = 1740.73
$ws.Range("E21").Value = "1740.73"

$ws.Range("A22").Value = "87811004_0222_JP"
$ws.Range("B22").Value = 32
$ws.Range("C22").Value = "JPY"
$ws.Range("D22").Value = 11340
$ws.Range("E22").Value = "11340"

$ws.Range("A23").Value = "87811004_0222_GB"
$ws.Range("B23").Value = 453
$ws.Range("C23").Value = "GBP"
$ws.Range("D23").Value = 1597.17
$ws.Range("E23").Value = "1597.17"

# Restore the default (General) style on column E so no residual text
# number-format is left applied to these cells, matching the target style.
$eRange.Style = "Normal"
